$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 150
$ws1.Range("F4").Value = 8
$ws1.Range("F9").Value = 195
$ws1.Range("G9").Value = 68
$ws1.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202408/Y6P7qrm11724139193256.jpeg"

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 150
$ws4.Range("F5").Value = 8
$ws4.Range("F10").Value = 195
$ws4.Range("G10").Value = 68
$ws4.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202408/Y6P7qrm11724139193256.jpeg"
